# EnemyTemplates.xlsx — "StatTable" sheet, EASY-mode STR (column G)
#
# Bring the EASY-mode STR-derivation formula up to the same ratio already
# used by the NORMAL-mode column (K): HEX2DEC($C)/9 -> HEX2DEC($C)*7/30
#
# The column is built from two shared-formula blocks (master cells G4 and
# G68, covering G4:G67 and G68:G112) plus one standalone formula in G3.
# Re-assigning .Formula across each full block keeps the shared-formula
# grouping intact (Excel/this engine re-derives the relative refs for every
# cell in the range from the single string we give it) and forces every
# cached <v> in the block to recalculate against the new ratio.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatTable")

$ws.Range("G3").Formula = '=IF($C3="00",$C3,DEC2HEX(MAX(ROUND(HEX2DEC($C3)*7/30,0),1),2))'
$ws.Range("G4:G67").Formula = '=IF($C4="00",$C4,DEC2HEX(MAX(ROUND(HEX2DEC($C4)*7/30,0),1),2))'
$ws.Range("G68:G112").Formula = '=IF($C68="00",$C68,DEC2HEX(MAX(ROUND(HEX2DEC($C68)*7/30,0),1),2))'

# Match the author's new selection — Glasses (H) and Ring (I) columns are now
# included alongside the STR (G) column for the Easy-mode block, anchored on
# F3 (the Easy-mode ATK column).
$ws.Range("F3:I112").Select()
